$d = $word.ActiveDocument

# In the "Institutional Environment"/training-opportunity paragraph, the
# standalone run containing just "Specifically" (immediately followed by
# " this project will provide training in splicing...") is missing the
# comma after the introductory adverb. Add it: "Specifically" -> "Specifically,"
#
# Find.Execute signature mirrors VBA:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,
#           ReplaceWith, Replace)
# wdFindContinue = 1 (Wrap), wdReplaceOne = 1 (Replace)
$find = $d.Content.Find
$found = $find.Execute("Specifically", $true, $true, $false, $false, $false, `
                        $true, 1, $false, "Specifically,", 1)

if (-not $found) {
    throw "Could not find the target run 'Specifically' to update."
}

Write-Output "Updated 'Specifically' -> 'Specifically,'"
